# Add season-record columns (Wins/Losses/Ties) to the roster/statistics
# sheet. Every player row gets the team's season record repeated across
# three new trailing columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column headers, styled like the existing header
# cells (bold font, centered/top-aligned, thin box border).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows (2-43): season record values, same for every row -
# 100 wins, 62 losses, 0 ties.
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 100
    $ws.Cells.Item($r, 31).Value = 62
    $ws.Cells.Item($r, 32).Value = 0
}
